$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 98; existing rows 98-105 shift down to 99-106
$ws.Rows.Item(98).Insert()

# Populate the newly inserted row 98 with the new price observation
$ws.Cells.Item(98, 1).Value = 8
$ws.Cells.Item(98, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(98, 3).Value = "Coquimbo"
$ws.Cells.Item(98, 4).Value = 44769
$ws.Cells.Item(98, 5).Value = 4
$ws.Cells.Item(98, 6).Value = 100112052
$ws.Cells.Item(98, 7).Value = "Albahaca"
$ws.Cells.Item(98, 8).Value = "Sin especificar"
$ws.Cells.Item(98, 9).Value = "Primera"
$ws.Cells.Item(98, 10).Value = 1400
$ws.Cells.Item(98, 11).Value = 3300
$ws.Cells.Item(98, 12).Value = 3500
$ws.Cells.Item(98, 13).Value = 3400
$ws.Cells.Item(98, 14).Value = "`$/paquete"
$ws.Cells.Item(98, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(98, 16).Value = 3400
$ws.Cells.Item(98, 17).Value = 1
$ws.Cells.Item(98, 18).Value = "Hortaliza"
